{"js": "// The document contains two test paragraphs (\"Hello World\" and\n// \"This is a new line\"). The commit removes that test content,\n// leaving the document body with a single empty paragraph.\nconst body = context.document.body;\nbody.clear();\nawait context.sync();\n", "ps1": "# The document contains two test paragraphs (\"Hello World\" and\n# \"This is a new line\"). The commit removes that test content,\n# leaving the document body with a single empty paragraph.\n$d = $word.ActiveDocument\n\n# Content.End sits just before the document's final paragraph mark, so\n# deleting that range collapses every paragraph but the last one down\n# to nothing; then we clear the remaining paragraph's text so the body\n# ends up with a single, empty paragraph.\nif ($d.Content.End -gt 0) {\n    $d.Range(0, $d.Content.End).Delete()\n}\n$d.Content.Text = \"\"\n"}
